$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the authoritative diff (cell ref -> new text value).
# Leading apostrophe forces Excel to store numeric-looking text as literal text
# (matches the original inlineStr/text cell type instead of being parsed as a number).
$updates = [ordered]@{
    "D2" = "'301.78"
    "E2" = "'0.76%"
    "D3" = "'32.90"
    "E3" = "'4.59%"
    "D4" = "'4.958"
    "E4" = "'-2.68%"
    "D5" = "'0.07780"
    "E5" = "'-1.32%"
    "D6" = "'1.968"
    "E6" = "'-13.02%"
    "D7" = "'7.852"
    "E7" = "'0.55%"
    "B8" = "MXToken"
    "C8" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "D8" = "'0.9274"
    "E8" = "'0.65%"
    "B9" = "WazirX"
    "C9" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "D9" = "'0.1768"
    "E9" = "'1.34%"
    "B10" = "LiechtensteinCryptoassetsExchange"
    "C10" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "D10" = "'0.07878"
    "E10" = "'3.79%"
    "B11" = "MandalaExchangeToken"
    "C11" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "D11" = "'0.08657"
    "E11" = "'-6.42%"
    "B12" = "BitrueCoin"
    "C12" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "D12" = "'0.03151"
    "E12" = "'4.81%"
    "B13" = "BitMartToken"
    "C13" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "D13" = "'0.1003"
    "E13" = "'0.05%"
    "B14" = "BitForexToken"
    "C14" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "D14" = "'0.001512"
    "E14" = "'-0.46%"
    "B15" = "TigerCash"
    "C15" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "D15" = "'0.005895"
    "E15" = "'-2.53%"
    "B16" = "LEO"
    "C16" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "D16" = "'3.462"
    "E16" = "'-0.45%"
    "B17" = "GateToken"
    "C17" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
    "D17" = "'3.799"
    "E17" = "'-1.34%"
    "D18" = "'2.155"
    "E18" = "'-3.95%"
    "E20" = "'0.84%"
    "D21" = "'4.329"
    "E21" = "'10.19%"
    "E22" = "'16.45%"
    "D23" = "'0.04558"
    "E23" = "'-1.15%"
    "E24" = "'-2.27%"
    "D25" = "'0.004437"
    "E25" = "'-0.84%"
    "E26" = "'0.12%"
    "D39" = "'0.01716"
    "E39" = "'-1.35%"
    "D40" = "'0.04720"
    "E40" = "'2.26%"
    "D41" = "'0.007884"
    "E41" = "'13.70%"
    "D42" = "'0.1354"
    "E42" = "'-0.43%"
    "E43" = "'10.28%"
    "D44" = "'0.01050"
    "E44" = "'1.96%"
    "D45" = "'0.00006251"
    "E45" = "'-0.50%"
    "E46" = "'0.11%"
    "B47" = "BOLO"
    "C47" = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
    "D47" = "'0.8234"
    "E47" = "'10.33%"
    "B48" = "CoinbaseStockToken"
    "C48" = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
    "D48" = "'0.003103"
    "E48" = "'-61.11%"
    "E49" = "'0.11%"
    "E50" = "'0.11%"
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

